$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.138.46"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +1.35%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.999.24"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +2.01%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.80"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.623"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.19%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.79"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.33%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("E9").Value = "  +2.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0804"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.61%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.103"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.58%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.97"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +5.20%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.44"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +5.26%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.294.34"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.04%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.840"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.42"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +2.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.000.12"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.10%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.073.09"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.11"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.50%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0862"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.32%  "
$ws.Range("E21").Value = "  +1.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "230.04"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("E24").Value = "  -0.21%  "
$ws.Range("E25").Value = "  -0.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.43"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.46%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.140"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -3.87%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "164.19"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +2.24%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.60"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.66%  "
$ws.Range("E30").Value = "  +12.11%  "
$ws.Range("E31").Value = "  +1.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.79"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.84%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0656"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +6.32%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.47"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.13%  "
$ws.Range("E35").Value = "  +4.20%  "
$ws.Range("E36").Value = "  +0.16%  "
$ws.Range("E37").Value = "  +2.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.31"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -7.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.36"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.68%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0981"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("E41").Value = "  +0.79%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0214"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.74%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.18"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.76%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.55"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +3.66%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "90.84"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +2.82%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.367.44"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.24%  "
$ws.Range("E47").Value = "  +1.30%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.35"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.93%  "
$ws.Range("E49").Value = "  +12.91%  "
$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "46.99"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +5.09%  "
$ws.Range("B51").Value = "MXToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.84"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.13%  "
